$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "Keda"
$ws.Name = "Keda"

# Row 6 ("Urban") - the whole data range becomes confidential-marker "…"
$ws.Range("C6:E6").Value = "…"
$ws.Range("G6").Value = "…"
$ws.Range("I6").Value = "…"

# Row 7 ("Rural") - three cells become the new "..." placeholder, others unchanged
$ws.Range("C7").Value = "..."
$ws.Range("E7").Value = "..."
$ws.Range("G7").Value = "..."

# Remove the blank row between the data table and the footnote row,
# shifting the footnote (row 9) up to row 8
$ws.Rows("8").Delete()
